$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.392.39'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.570.26'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3743'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.81%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.47'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3365'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07465'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.129'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.914'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.864'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.21%  '
$ws.Range("D16").Value = '1.569.20'
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("E17").Value = '  -2.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.54%  '
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.164'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.86'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.14%  '
$ws.Range("D24").Value = '22.384.56'
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.365'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.547'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.988'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("D31").Value = '1.745.42'
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9956'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.952'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.912'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.716'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08421'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.381'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02451'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.70%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2246'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.51%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06460'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.369'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.53%  '
$ws.Range("E42").Value = '  -3.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6200'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.811'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5788'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.62%  '
$ws.Range("E48").Value = '  -2.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("E50").Value = '  -6.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07297'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.11%  '
